$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of dish data appended below the existing table (rows 7-13).
# Columns: A = DishName, C = DishUrl (B left blank, matching the source diff)
$ws.Range("C7").Value = "https://res.cloudinary.com/rainforest-cruises/images/c_fill,g_auto/f_auto,q_auto/v1661887113/indian-food/indian-food.jpg"

$ws.Range("A8").Value = "Chole Bhature"
$ws.Range("C8").Value = "https://www.bitesbee.com/wp-content/uploads/2021/09/Chole-Bhature.jpg"

$ws.Range("C9").Value = "https://static.toiimg.com/photo/96559646/96559646.jpg"

$ws.Range("A10").Value = "Chiken Pokoda"
$ws.Range("C10").Value = "https://res.cloudinary.com/rainforest-cruises/images/c_fill,g_auto/f_auto,q_auto/w_1120,h_732,c_fill,g_auto/v1661347434/india-food-chicken-pakora/india-food-chicken-pakora-1120x732.jpg"

$ws.Range("A11").Value = "Misal Pav"
$ws.Range("C11").Value = "https://www.holidify.com/blog/wp-content/uploads/2015/11/Maharashtras_Misal_Pav.jpg"

$ws.Range("C12").Value = "https://www.eatingwell.com/thmb/kApzxiwARsrq9ILeAkRRdvlXrQs=/1500x0/filters:no_upscale():max_bytes(150000):strip_icc()/3879366-0fbeae20516648df942721efa761894c.jpg"

$ws.Range("A13").Value = "Biryani"
$ws.Range("C13").Value = "https://images.yummy.ph/yummy/uploads/2023/02/biryani.jpg"

# Update the view's active cell/selection to match the final saved state.
$ws.Range("A8").Select()
